# "finished up the project files"
# - Updates timing numbers for the three existing search blocks.
# - Fixes a typo'd result count (C17: 18151 -> 18153).
# - Renames the "astar_search" header to "astar_search with h_1" and moves
#   it (with its data block) further down the sheet.
# - Adds two brand-new result blocks for h_ignore_preconditions and
#   h_pg_levelsum heuristics.
# - Updates window/view state to where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Section numbers in column A (only on the header row of each block)
# ---------------------------------------------------------------------
$ws.Range("A1").Value2  = 1
$ws.Range("A7").Value2  = 3
$ws.Range("A13").Value2 = 5

# ---------------------------------------------------------------------
# Refresh the timing ("Time elapsed") values + the one corrected count
# for the three blocks that already existed (breadth_first_search,
# depth_first_graph_search, astar_search/h_1).
# ---------------------------------------------------------------------
$ws.Range("E3").Value2  = 0.078426908992696498
$ws.Range("E4").Value2  = 28.952266854001198
$ws.Range("E5").Value2  = 182.17717483200201

$ws.Range("E9").Value2  = 0.026630198000930201
$ws.Range("E10").Value2 = 7.2162615499982996
$ws.Range("E11").Value2 = 3.1420503469998899

$ws.Range("E15").Value2 = 0.068847660993924295
$ws.Range("E16").Value2 = 22.126883590972199
$ws.Range("C17").Value2 = 18153
$ws.Range("E17").Value2 = 115.730096187005

# ---------------------------------------------------------------------
# The old row-13 block's header ("astar_search") is blanked out; the
# heading text is reborn further down as "astar_search with h_1".
# ---------------------------------------------------------------------
$ws.Range("B13").Value2 = ""

# ---------------------------------------------------------------------
# New block: astar_search with h_1 (rows 21-25) - reuses the shared
# string that used to say plain "astar_search".
# ---------------------------------------------------------------------
$ws.Range("A21").Value2 = 8
$ws.Range("B21").Value2 = "astar_search with h_1"
$ws.Range("B21:F21").Merge()
$ws.Range("B21:F21").HorizontalAlignment = -4108

$ws.Range("B22").Value2 = "Expansions"
$ws.Range("C22").Value2 = "Goal Tests"
$ws.Range("D22").Value2 = "New Nodes"
$ws.Range("E22").Value2 = "Time elapsed"
$ws.Range("F22").Value2 = "length"

$ws.Range("A23").Value2 = "Problem 1"
$ws.Range("B23").Value2 = 55
$ws.Range("C23").Value2 = 57
$ws.Range("D23").Value2 = 224
$ws.Range("E23").Value2 = 0.085609160014428198
$ws.Range("F23").Value2 = 6

$ws.Range("A24").Value2 = "Problem 2"
$ws.Range("B24").Value2 = 4853
$ws.Range("C24").Value2 = 4855
$ws.Range("D24").Value2 = 44041
$ws.Range("E24").Value2 = 20.296631225995899
$ws.Range("F24").Value2 = 9

$ws.Range("A25").Value2 = "Problem 3"
$ws.Range("B25").Value2 = 18151
$ws.Range("C25").Value2 = 18153
$ws.Range("D25").Value2 = 159038
$ws.Range("E25").Value2 = 100.67678316100501
$ws.Range("F25").Value2 = 12

# ---------------------------------------------------------------------
# New block: astar_search with h_ignore_preconditions (rows 29-33)
# ---------------------------------------------------------------------
$ws.Range("A29").Value2 = 9
$ws.Range("B29").Value2 = "astar_search with h_ignore_preconditions"
$ws.Range("B29:F29").Merge()
$ws.Range("B29:F29").HorizontalAlignment = -4108

$ws.Range("B30").Value2 = "Expansions"
$ws.Range("C30").Value2 = "Goal Tests"
$ws.Range("D30").Value2 = "New Nodes"
$ws.Range("E30").Value2 = "Time elapsed"
$ws.Range("F30").Value2 = "length"

$ws.Range("A31").Value2 = "Problem 1"
$ws.Range("B31").Value2 = 41
$ws.Range("C31").Value2 = 43
$ws.Range("D31").Value2 = 170
$ws.Range("E31").Value2 = 0.067934118007542496
$ws.Range("F31").Value2 = 6

$ws.Range("A32").Value2 = "Problem 2"
$ws.Range("B32").Value2 = 1450
$ws.Range("C32").Value2 = 1452
$ws.Range("D32").Value2 = 13303
$ws.Range("E32").Value2 = 7.9578909480187496
$ws.Range("F32").Value2 = 9

$ws.Range("A33").Value2 = "Problem 3"
$ws.Range("B33").Value2 = 5038
$ws.Range("C33").Value2 = 5040
$ws.Range("D33").Value2 = 44926
$ws.Range("E33").Value2 = 32.798792346991803
$ws.Range("F33").Value2 = 12

# ---------------------------------------------------------------------
# New block: astar_search with h_pg_levelsum (rows 36-40) - Problem 3
# run hadn't finished yet, so only the row labels are there for
# Problem 1/2, with the completed Problem 3 numbers below.
# ---------------------------------------------------------------------
$ws.Range("A36").Value2 = 10
$ws.Range("B36").Value2 = "astar_search with h_pg_levelsum"
$ws.Range("B36:F36").Merge()
$ws.Range("B36:F36").HorizontalAlignment = -4108

$ws.Range("B37").Value2 = "Expansions"
$ws.Range("C37").Value2 = "Goal Tests"
$ws.Range("D37").Value2 = "New Nodes"
$ws.Range("E37").Value2 = "Time elapsed"
$ws.Range("F37").Value2 = "length"

$ws.Range("A38").Value2 = "Problem 1"
$ws.Range("A39").Value2 = "Problem 2"

$ws.Range("A40").Value2 = "Problem 3"
$ws.Range("B40").Value2 = 0
$ws.Range("C40").Value2 = 2
$ws.Range("D40").Value2 = 0
$ws.Range("E40").Value2 = 0.31762763799633797
$ws.Range("F40").Value2 = 0

# ---------------------------------------------------------------------
# View state: scroll down to the newest block and leave the selection
# where the author last left it.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$win.Left = 13120
$win.Top = 0
$ws.Range("F32").Select()
